# Rename the worksheet from "Property1" to "DataNode" to unify the
# DataNode / DataTable / Entity naming convention used across the
# config workbooks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Update the frozen-pane selection to the range that was selected when
# the file was last saved (A9:N35).
$ws.Range("A9:N35").Select()
